$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Core input edits ---
# Sheet1: Price assumption changes from 81 to 76.3, and the cell's number
# format is widened to show two decimals.
$ws1.Range("E2").Value = 76.3
$ws1.Range("E2").NumberFormat = "#,##0.00"

# Sheet2: ROIC/reinvestment-rate assumption (Z103) and the long-run decline
# rate (Z104) are both lowered.
$ws2.Range("Z103").Value = 0.04
$ws2.Range("Z104").Value = -0.01

# Z107 (price target) also gets a two-decimal number format applied.
$ws2.Range("Z107").NumberFormat = "#,##0.00"

$excel.CalculateFullRebuild()

# --- View / selection state ---
# Leave a selection of E2 behind on Sheet1 without making it the active tab.
$ws1.Range("E2").Select()

# Sheet2 stays the active/visible tab, with its frozen pane scrolled down and
# the bottom-right pane's selection moved to Z106.
$ws2.Activate()
$ws2.Range("M85").Select()
$excel.ActiveWindow.ScrollRow = 85
$excel.ActiveWindow.ScrollColumn = 13
$ws2.Range("Z106").Select()

# --- Print setup ---
$ws1.PageSetup.Orientation = 1
